# Rewrite the M2Doc field "{m:'dh1.gif'.asImage().fit(200, 500)}" from a
# real Word field (fldChar begin/instrText.../fldChar end) into plain
# literal text runs ({ ... }), matching the TokenIteratorFieldRewriterSplit
# output format, while preserving the run-level colour formatting used on
# the statement body.

$d = $word.ActiveDocument

# Locate the (only) field in the document and the paragraph that contains it.
$field = $d.Fields(1)
$para = $field.Code.Paragraphs(1)
$rng = $para.Range

# Build the replacement paragraph content as WordprocessingML, wrapped in the
# pkg:package envelope Range.InsertXML / Range.WordOpenXML use. Each former
# <w:instrText> run becomes a <w:t> run with identical text/rPr; the leading
# "fldChar begin" + space become a single "{" run, and the trailing space +
# "fldChar end" become a single "}" run (preserving xml:space="preserve").
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r w:rsidR="004B598D"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>'dh</w:t></w:r><w:r w:rsidR="00321AA1"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>1</w:t></w:r><w:r w:rsidR="004B598D"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>.gif'.asImage()</w:t></w:r><w:r w:rsidR="0047710F"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t xml:space="preserve">.fit(200, </w:t></w:r><w:r w:rsidR="00DD079F"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>5</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="0047710F"><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>00)</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$rng.InsertXML($xml)
